$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.430437207221985
$ws.Range("B1").Value = 3.761228799819946
$ws.Range("C1").Value = 3.48394250869751
$ws.Range("D1").Value = 3.343560457229614
$ws.Range("E1").Value = 1.207168459892273
